# Insert a new weekly record as row 251 (pushing the existing rows 251..287
# down to 252..288) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 251, shifting everything
# below it (including the former row 251) down by one.
$ws.Rows.Item(251).Insert()

# Populate the freshly inserted row 251 with the new data point.
$ws.Cells.Item(251, 1).Value  = 6
$ws.Cells.Item(251, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(251, 3).Value  = "Metropolitana"
$ws.Cells.Item(251, 4).Value  = 44474
$ws.Cells.Item(251, 5).Value  = 13
$ws.Cells.Item(251, 6).Value  = 100112039
$ws.Cells.Item(251, 7).Value  = "Ciboulette"
$ws.Cells.Item(251, 8).Value  = "Sin especificar"
$ws.Cells.Item(251, 9).Value  = "Primera"
$ws.Cells.Item(251, 10).Value = 850
$ws.Cells.Item(251, 11).Value = 800
$ws.Cells.Item(251, 12).Value = 1000
$ws.Cells.Item(251, 13).Value = 892
$ws.Cells.Item(251, 14).Value = "`$/docena de atados"
$ws.Cells.Item(251, 15).Value = "Región Metropolitana"
$ws.Cells.Item(251, 16).Value = 297
$ws.Cells.Item(251, 17).Value = 3
$ws.Cells.Item(251, 18).Value = "Hortaliza"
